# Add the new "remove kth node from end of linked list" row (row 3) to the
# "链表" (linked list) problem-list sheet, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 data: No. / lc / 题目 / 解题方法 / 解题关键词 / 时间复杂度 / 空间复杂度
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = "给定一个链表，删除链表的倒数第n个节点，并且返回链表的头结点"
$ws.Range("D3").Value = "`n1. 快指针fast先从solder向前走n步`n2. 判断fast为空，就停止；不为空就继续向下`n3. fast没有达到尾部，fast与slow指针同时向前走一 步，slow走第一步时是从solder走的`n4. fast到达链表尾部(非空节点)，slow就指向倒数第n个节点"
$ws.Range("E3").Value = "快慢指针`n滑动窗口`nsolder"
$ws.Range("F3").Value = "O(L)，L是窗口长度"
$ws.Range("G3").Value = "O(1)"

# The new row holds multi-line wrapped text, so it needs a taller row (140pt),
# matching the other wrapped-text data row.
$ws.Rows.Item(3).RowHeight = 140

# The author's last action before saving was selecting C3.
$ws.Range("C3").Select()
